# The site footer block ("Ver no Jupiter ..." / copyright line), plus the
# blank paragraph separating it from the "LOT2041: ..." requisito line,
# was dropped from this page in the rebuild. Find the paragraph holding
# the requisito text and the paragraph holding the copyright notice, then
# delete everything from just after the former through the end of the
# latter (this removes the blank paragraph, the "Ver no Jupiter ..."
# paragraph, and the copyright paragraph, while leaving the following
# blank paragraph and the page-break paragraph untouched).

$d = $word.ActiveDocument
$paragraphs = $d.Paragraphs
$n = $paragraphs.Count

$reqIndex = 0
$copyrightIndex = 0
for ($i = 1; $i -le $n; $i++) {
    $t = $paragraphs.Item($i).Range.Text
    if ($t -like "LOT2041: Engenharia Bioqu*mica (Requisito fraco)*") {
        $reqIndex = $i
    }
    if ($t -like "*2020 . Contact: luizeleno@usp.br*") {
        $copyrightIndex = $i
    }
}

$from = $paragraphs.Item($reqIndex + 1)
$to = $paragraphs.Item($copyrightIndex)
$range = $d.Range($from.Range.Start, $to.Range.End)
$range.Delete()
